$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, matching the formatting of the existing
# header cells (e.g. G1 -- bold, bordered, centered) by copying its format.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add the data value for the new "Save" column in row 2
$ws.Range("H2").Value = 1
